$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record as row 474 (pushing the existing rows 474:483 down to 475:484)
$ws.Rows("474:474").Insert()

# Populate the newly inserted row with the new weekly price observation
$ws.Range("A474").Value = 4
$ws.Range("B474").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C474").Value = "Los Lagos"
$ws.Range("D474").Value = 45239
$ws.Range("E474").Value = 10
$ws.Range("F474").Value = 100112021
$ws.Range("G474").Value = "Ají"
$ws.Range("H474").Value = "Inferno"
$ws.Range("I474").Value = "Primera"
$ws.Range("J474").Value = 70
$ws.Range("K474").Value = 42000
$ws.Range("L474").Value = 42000
$ws.Range("M474").Value = 42000
$ws.Range("N474").Value = "$/caja 10 kilos"
$ws.Range("O474").Value = "Región de Arica y Parinacota"
$ws.Range("P474").Value = 4200
$ws.Range("Q474").Value = 10
$ws.Range("R474").Value = "Hortaliza"
